$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 15083.375
$ws.Range("I40").Value = 18913
$ws.Range("K40").Value = 18913
$ws.Range("M40").Value = -18738
$ws.Range("H111").Value = 13895550
$ws.Range("I111").Value = 15631622
$ws.Range("K111").Value = 46894866
$ws.Range("M111").Value = -46891799
$ws.Range("H118").Value = 3499.889
$ws.Range("I118").Value = 3826.125
$ws.Range("K118").Value = 11478.375
$ws.Range("M118").Value = -9821.375
$ws.Range("H129").Value = 1872.25
$ws.Range("I129").Value = 1872.25
$ws.Range("K129").Value = 5616.75
$ws.Range("M129").Value = -616.75
$ws.Range("H135").Value = 357697.06
$ws.Range("I135").Value = 435158.25
$ws.Range("J135").Value = 1375.6
$ws.Range("K135").Value = 3916424.25
$ws.Range("L135").Value = 12380.4
$ws.Range("M135").Value = -3913889.25
$ws.Range("N135").Value = -17450.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 47744.77
$ws.Range("I74").Value = 58422.82
$ws.Range("K74").Value = 58422.82
$ws.Range("M74").Value = -57548.82
$ws.Range("H77").Value = 47744.77
$ws.Range("I77").Value = 58422.82
$ws.Range("K77").Value = 292114.1
$ws.Range("M77").Value = -287746.1

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 50003816
$ws.Range("I86").Value = 2687
$ws.Range("K86").Value = 2687
$ws.Range("M86").Value = -1564
$ws.Range("H89").Value = 50003816
$ws.Range("I89").Value = 2687
$ws.Range("K89").Value = 13435
$ws.Range("M89").Value = -7819
$ws.Range("H99").Value = 11365845
$ws.Range("I99").Value = 1005
$ws.Range("K99").Value = 1005
$ws.Range("M99").Value = 493
$ws.Range("H105").Value = 2714.9429
$ws.Range("I105").Value = 2393.3215
$ws.Range("J105").Value = 4001.4285
$ws.Range("K105").Value = 2393.3215
$ws.Range("L105").Value = 4001.4285
$ws.Range("M105").Value = -646.3215
$ws.Range("N105").Value = -7495.4285
$ws.Range("H134").Value = 3727.423
$ws.Range("I134").Value = 2016.3125
$ws.Range("K134").Value = 6048.9375
$ws.Range("M134").Value = -3513.9375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4103
$ws.Range("I16").Value = 2098.3
$ws.Range("J16").Value = 5534.9287
$ws.Range("K16").Value = 2098.3
$ws.Range("L16").Value = 5534.9287
$ws.Range("M16").Value = -1811.3
$ws.Range("N16").Value = -6108.9287
$ws.Range("H31").Value = 6690.4204
$ws.Range("I31").Value = 2918.5789
$ws.Range("J31").Value = 11313.968
$ws.Range("K31").Value = 2918.5789
$ws.Range("L31").Value = 11313.968
$ws.Range("M31").Value = -2623.5789
$ws.Range("N31").Value = -11903.968
$ws.Range("H34").Value = 6690.4204
$ws.Range("I34").Value = 2918.5789
$ws.Range("J34").Value = 11313.968
$ws.Range("K34").Value = 2918.5789
$ws.Range("L34").Value = 11313.968
$ws.Range("M34").Value = -2716.5789
$ws.Range("N34").Value = -11717.968
$ws.Range("H86").Value = 3910640
$ws.Range("I86").Value = 5212521
$ws.Range("K86").Value = 5212521
$ws.Range("M86").Value = -5211398
$ws.Range("H89").Value = 3910640
$ws.Range("I89").Value = 5212521
$ws.Range("K89").Value = 26062605
$ws.Range("M89").Value = -26056989
$ws.Range("H113").Value = 4103
$ws.Range("I113").Value = 2098.3
$ws.Range("J113").Value = 5534.9287
$ws.Range("K113").Value = 2098.3
$ws.Range("L113").Value = 5534.9287
$ws.Range("M113").Value = 71.69999999999982
$ws.Range("N113").Value = -9874.9287

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 2000670.2
$ws.Range("J12").Value = 3125565
$ws.Range("L12").Value = 9376695
$ws.Range("N12").Value = -9377041
$ws.Range("H34").Value = 3610.8572
$ws.Range("I34").Value = 339.5
$ws.Range("J34").Value = 5624
$ws.Range("K34").Value = 1018.5
$ws.Range("L34").Value = 16872
$ws.Range("M34").Value = -934.5
$ws.Range("N34").Value = -17040
$ws.Range("H136").Value = 2009.6666
$ws.Range("I136").Value = 2009.6666
$ws.Range("K136").Value = 6028.9998
$ws.Range("M136").Value = -928.9997999999996
$ws.Range("H137").Value = 200917.6
$ws.Range("I137").Value = 125784.625
$ws.Range("K137").Value = 377353.875
$ws.Range("M137").Value = -372253.875
$ws.Range("H138").Value = 3781.2307
$ws.Range("I138").Value = 3696
$ws.Range("K138").Value = 11088
$ws.Range("M138").Value = -5948

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5671.98
$ws.Range("I70").Value = 4474
$ws.Range("K70").Value = 4474
$ws.Range("M70").Value = -4204
$ws.Range("H73").Value = 5671.98
$ws.Range("I73").Value = 4474
$ws.Range("K73").Value = 4474
$ws.Range("M73").Value = -3538
$ws.Range("H97").Value = 2431.6785
$ws.Range("I97").Value = 2135.1765
$ws.Range("J97").Value = 2889.9092
$ws.Range("K97").Value = 2135.1765
$ws.Range("L97").Value = 2889.9092
$ws.Range("M97").Value = -1639.1765
$ws.Range("N97").Value = -3881.9092
$ws.Range("H113").Value = 8164.2383
$ws.Range("I113").Value = 4689.8
$ws.Range("J113").Value = 9250
$ws.Range("K113").Value = 4689.8
$ws.Range("L113").Value = 9250
$ws.Range("M113").Value = -2519.8
$ws.Range("N113").Value = -13590
$ws.Range("H126").Value = 3009.3333
$ws.Range("I126").Value = 2896
$ws.Range("J126").Value = 3041.7144
$ws.Range("K126").Value = 8688
$ws.Range("L126").Value = 9125.143199999999
$ws.Range("M126").Value = -6218
$ws.Range("N126").Value = -14065.1432
$ws.Range("H134").Value = 99995.336
$ws.Range("J134").Value = 99995.336
$ws.Range("L134").Value = 299986.008
$ws.Range("N134").Value = -305056.008

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6223.8096
$ws.Range("I7").Value = 4739.3
$ws.Range("J7").Value = 7573.364
$ws.Range("K7").Value = 4739.3
$ws.Range("L7").Value = 7573.364
$ws.Range("M7").Value = -4627.3
$ws.Range("N7").Value = -7797.364
$ws.Range("H46").Value = 29120534
$ws.Range("I46").Value = 34482760
$ws.Range("J46").Value = 27779978
$ws.Range("K46").Value = 34482760
$ws.Range("L46").Value = 27779978
$ws.Range("M46").Value = -34482572
$ws.Range("N46").Value = -27780354
$ws.Range("H61").Value = 5980.8823
$ws.Range("I61").Value = 5253.143
$ws.Range("J61").Value = 6490.3
$ws.Range("K61").Value = 5253.143
$ws.Range("L61").Value = 6490.3
$ws.Range("M61").Value = -5051.143
$ws.Range("N61").Value = -6894.3
$ws.Range("H68").Value = 2999.5
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 2999.5
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 2999.5
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -4497.5
$ws.Range("H71").Value = 2999.5
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 2999.5
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 14997.5
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -22485.5
$ws.Range("H113").Value = 5980.8823
$ws.Range("I113").Value = 5253.143
$ws.Range("J113").Value = 6490.3
$ws.Range("K113").Value = 5253.143
$ws.Range("L113").Value = 6490.3
$ws.Range("M113").Value = -3083.143
$ws.Range("N113").Value = -10830.3
$ws.Range("H126").Value = 6223.8096
$ws.Range("I126").Value = 4739.3
$ws.Range("J126").Value = 7573.364
$ws.Range("K126").Value = 14217.9
$ws.Range("L126").Value = 22720.092
$ws.Range("M126").Value = -11747.9
$ws.Range("N126").Value = -27660.092

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 20000
$ws.Range("I62").Value = 20000
$ws.Range("K62").Value = 20000
$ws.Range("M62").Value = -19376
$ws.Range("H65").Value = 20000
$ws.Range("I65").Value = 20000
$ws.Range("K65").Value = 100000
$ws.Range("M65").Value = -96880
$ws.Range("H81").Value = 18194134
$ws.Range("I81").Value = 1639.1428
$ws.Range("K81").Value = 3278.2856
$ws.Range("M81").Value = -2217.2856
$ws.Range("H84").Value = 18194134
$ws.Range("I84").Value = 1639.1428
$ws.Range("K84").Value = 16391.428
$ws.Range("M84").Value = -11087.428
